$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 'sv'
$ws.Range("J3").Value = 'Statement-opinion'
$ws.Range("I8").Value = 'sd'
$ws.Range("J8").Value = 'Statement-non-opinion'
$ws.Range("I11").Value = 'sv'
$ws.Range("J11").Value = 'Statement-opinion'
$ws.Range("I17").Value = 'sv'
$ws.Range("J17").Value = 'Statement-opinion'
$ws.Range("I19").Value = 'sd'
$ws.Range("J19").Value = 'Statement-non-opinion'
$ws.Range("I21").Value = 'sd'
$ws.Range("J21").Value = 'Statement-non-opinion'
$ws.Range("I23").Value = 'qy'
$ws.Range("J23").Value = 'Yes-No-Question'
$ws.Range("I26").Value = 'sd'
$ws.Range("J26").Value = 'Statement-non-opinion'
$ws.Range("I36").Value = 'sd'
$ws.Range("J36").Value = 'Statement-non-opinion'
$ws.Range("I39").Value = 'ba'
$ws.Range("J39").Value = 'Appreciation'
$ws.Range("I43").Value = 'sd'
$ws.Range("J43").Value = 'Statement-non-opinion'
$ws.Range("I48").Value = 'ba'
$ws.Range("J48").Value = 'Appreciation'
$ws.Range("I51").Value = 'sd'
$ws.Range("J51").Value = 'Statement-non-opinion'
$ws.Range("I71").Value = 'sv'
$ws.Range("J71").Value = 'Statement-opinion'
$ws.Range("I84").Value = 'b'
$ws.Range("J84").Value = 'Acknowledge (Backchannel)'
$ws.Range("I96").Value = 'aa'
$ws.Range("J96").Value = 'Agree/Accept'
$ws.Range("I99").Value = 'aa'
$ws.Range("J99").Value = 'Agree/Accept'
$ws.Range("I100").Value = 'aa'
$ws.Range("J100").Value = 'Agree/Accept'
$ws.Range("I102").Value = 'sd'
$ws.Range("J102").Value = 'Statement-non-opinion'
$ws.Range("I103").Value = 'ba'
$ws.Range("J103").Value = 'Appreciation'
$ws.Range("I120").Value = 'sd'
$ws.Range("J120").Value = 'Statement-non-opinion'
$ws.Range("I131").Value = 'sd'
$ws.Range("J131").Value = 'Statement-non-opinion'
$ws.Range("I135").Value = 'sd'
$ws.Range("J135").Value = 'Statement-non-opinion'
$ws.Range("I140").Value = 'sd'
$ws.Range("J140").Value = 'Statement-non-opinion'
$ws.Range("I142").Value = 'sd'
$ws.Range("J142").Value = 'Statement-non-opinion'
$ws.Range("I146").Value = 'aa'
$ws.Range("J146").Value = 'Agree/Accept'
$ws.Range("I149").Value = 'ba'
$ws.Range("J149").Value = 'Appreciation'
$ws.Range("I156").Value = 'sv'
$ws.Range("J156").Value = 'Statement-opinion'
$ws.Range("I157").Value = 'sv'
$ws.Range("J157").Value = 'Statement-opinion'
$ws.Range("I158").Value = 'sv'
$ws.Range("J158").Value = 'Statement-opinion'
$ws.Range("I177").Value = 'sd'
$ws.Range("J177").Value = 'Statement-non-opinion'
$ws.Range("I183").Value = 'aa'
$ws.Range("J183").Value = 'Agree/Accept'
$ws.Range("I184").Value = 'aa'
$ws.Range("J184").Value = 'Agree/Accept'
$ws.Range("I191").Value = 'sd'
$ws.Range("J191").Value = 'Statement-non-opinion'
$ws.Range("I209").Value = 'sd'
$ws.Range("J209").Value = 'Statement-non-opinion'
$ws.Range("I210").Value = '%'
$ws.Range("J210").Value = 'Uninterpretable'
$ws.Range("I217").Value = 'sv'
$ws.Range("J217").Value = 'Statement-opinion'
$ws.Range("I220").Value = 'sv'
$ws.Range("J220").Value = 'Statement-opinion'
$ws.Range("I222").Value = '%'
$ws.Range("J222").Value = 'Uninterpretable'
$ws.Range("I241").Value = 'b'
$ws.Range("J241").Value = 'Acknowledge (Backchannel)'
$ws.Range("I255").Value = 'sv'
$ws.Range("J255").Value = 'Statement-opinion'
$ws.Range("I256").Value = 'sd'
$ws.Range("J256").Value = 'Statement-non-opinion'
$ws.Range("I257").Value = '%'
$ws.Range("J257").Value = 'Uninterpretable'
$ws.Range("I272").Value = 'ba'
$ws.Range("J272").Value = 'Appreciation'
$ws.Range("I288").Value = 'sv'
$ws.Range("J288").Value = 'Statement-opinion'
$ws.Range("I292").Value = 'sd'
$ws.Range("J292").Value = 'Statement-non-opinion'
$ws.Range("I293").Value = 'ba'
$ws.Range("J293").Value = 'Appreciation'
$ws.Range("I297").Value = '%'
$ws.Range("J297").Value = 'Uninterpretable'
$ws.Range("I300").Value = 'sd'
$ws.Range("J300").Value = 'Statement-non-opinion'
$ws.Range("I318").Value = 'sv'
$ws.Range("J318").Value = 'Statement-opinion'
$ws.Range("I330").Value = 'sd'
$ws.Range("J330").Value = 'Statement-non-opinion'
$ws.Range("I333").Value = 'sv'
$ws.Range("J333").Value = 'Statement-opinion'
$ws.Range("I339").Value = 'sv'
$ws.Range("J339").Value = 'Statement-opinion'
$ws.Range("I355").Value = 'sv'
$ws.Range("J355").Value = 'Statement-opinion'
$ws.Range("I357").Value = 'sv'
$ws.Range("J357").Value = 'Statement-opinion'
$ws.Range("I380").Value = 'sd'
$ws.Range("J380").Value = 'Statement-non-opinion'
$ws.Range("I385").Value = 'sv'
$ws.Range("J385").Value = 'Statement-opinion'
$ws.Range("I417").Value = 'aa'
$ws.Range("J417").Value = 'Agree/Accept'
$ws.Range("I426").Value = 'sv'
$ws.Range("J426").Value = 'Statement-opinion'
$ws.Range("I432").Value = 'sd'
$ws.Range("J432").Value = 'Statement-non-opinion'
$ws.Range("I442").Value = 'sv'
$ws.Range("J442").Value = 'Statement-opinion'
$ws.Range("I450").Value = 'aa'
$ws.Range("J450").Value = 'Agree/Accept'
$ws.Range("I453").Value = 'b'
$ws.Range("J453").Value = 'Acknowledge (Backchannel)'
$ws.Range("I456").Value = 'b'
$ws.Range("J456").Value = 'Acknowledge (Backchannel)'
$ws.Range("I485").Value = 'sv'
$ws.Range("J485").Value = 'Statement-opinion'
$ws.Range("I499").Value = 'sd'
$ws.Range("J499").Value = 'Statement-non-opinion'
$ws.Range("I513").Value = 'b'
$ws.Range("J513").Value = 'Acknowledge (Backchannel)'
$ws.Range("I514").Value = 'sd'
$ws.Range("J514").Value = 'Statement-non-opinion'
$ws.Range("I524").Value = '%'
$ws.Range("J524").Value = 'Uninterpretable'
$ws.Range("I532").Value = 'aa'
$ws.Range("J532").Value = 'Agree/Accept'
$ws.Range("I539").Value = '%'
$ws.Range("J539").Value = 'Uninterpretable'
$ws.Range("I547").Value = 'sv'
$ws.Range("J547").Value = 'Statement-opinion'
$ws.Range("I549").Value = 'sv'
$ws.Range("J549").Value = 'Statement-opinion'
$ws.Range("I553").Value = '%'
$ws.Range("J553").Value = 'Uninterpretable'
$ws.Range("I555").Value = 'sv'
$ws.Range("J555").Value = 'Statement-opinion'
$ws.Range("I557").Value = 'sd'
$ws.Range("J557").Value = 'Statement-non-opinion'
$ws.Range("I577").Value = 'sv'
$ws.Range("J577").Value = 'Statement-opinion'
$ws.Range("I579").Value = 'sd'
$ws.Range("J579").Value = 'Statement-non-opinion'
$ws.Range("I580").Value = 'sv'
$ws.Range("J580").Value = 'Statement-opinion'
$ws.Range("I581").Value = 'sd'
$ws.Range("J581").Value = 'Statement-non-opinion'
$ws.Range("I588").Value = 'sd'
$ws.Range("J588").Value = 'Statement-non-opinion'
$ws.Range("I591").Value = 'ba'
$ws.Range("J591").Value = 'Appreciation'
$ws.Range("I596").Value = 'sv'
$ws.Range("J596").Value = 'Statement-opinion'
$ws.Range("I597").Value = 'sv'
$ws.Range("J597").Value = 'Statement-opinion'
$ws.Range("I601").Value = 'sv'
$ws.Range("J601").Value = 'Statement-opinion'
$ws.Range("I605").Value = 'aa'
$ws.Range("J605").Value = 'Agree/Accept'
$ws.Range("I632").Value = 'ba'
$ws.Range("J632").Value = 'Appreciation'
$ws.Range("I643").Value = 'sd'
$ws.Range("J643").Value = 'Statement-non-opinion'
$ws.Range("I652").Value = 'aa'
$ws.Range("J652").Value = 'Agree/Accept'
$ws.Range("I653").Value = 'aa'
$ws.Range("J653").Value = 'Agree/Accept'
$ws.Range("I675").Value = 'aa'
$ws.Range("J675").Value = 'Agree/Accept'
